# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

function Repeat-Str($ch, $n) {
    $result = ""
    for ($i = 0; $i -lt $n; $i++) {
        $result = $result + $ch
    }
    return $result
}

# Force a TextRange (typically obtained via .Characters(start, len)) to take
# on new text as a *single run* with the formatting of the range's current
# first run. If the new text happens to match the current text exactly the
# backing store treats it as a no-op and keeps existing run boundaries, so
# we first scramble the range (same length, so the range stays valid) and
# then apply the real text.
function Set-RunText($range, $text) {
    $len = $range.Length
    if ($len -gt 0) {
        $placeholder = Repeat-Str "*" $len
        $range.Text = $placeholder
    }
    $range.Text = $text
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master + notes master date fields: 9/29/2014 -> 10/13/2014
# ---------------------------------------------------------------------------

$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $shp = $hm.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "9/29/2014") {
            $tr.Text = "10/13/2014"
        }
    }
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $shp = $nm.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "9/29/2014") {
            $tr.Text = "10/13/2014"
        }
    }
}

Write-Output "Step 1 done: date fields updated"
